$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 37.98277566666666
$ws.Range("H2").Value = 113.948327
$ws.Range("I2").Value = 0.697850645410475
$ws.Range("J2").Value = 0.6978506454104751
$ws.Range("M2").Value = 8.540560666666666
$ws.Range("N2").Value = 25.621682
$ws.Range("O2").Value = 0.4159358086620884
$ws.Range("P2").Value = 0.4159358086620884
$ws.Range("Q2").Value = 324.394199869557
$ws.Range("R2").Value = 2919.547798826014
$ws.Range("S2").Value = 0.2902610725241662
$ws.Range("T2").Value = 0.2902610725241663
$ws.Range("G3").Value = 37.98277566666666
$ws.Range("H3").Value = 113.948327
$ws.Range("I3").Value = 0.697850645410475
$ws.Range("J3").Value = 0.6978506454104751
$ws.Range("O3").Value = 0.563694901924408
$ws.Range("P3").Value = 0.563694901924408
$ws.Range("Q3").Value = 439.6335993972427
$ws.Range("R3").Value = 3956.702394575185
$ws.Range("S3").Value = 0.3933748511225425
$ws.Range("T3").Value = 0.3933748511225426
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("G4").Value = 37.98277566666666
$ws.Range("H4").Value = 113.948327
$ws.Range("I4").Value = 0.697850645410475
$ws.Range("J4").Value = 0.6978506454104751
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.026642
$ws.Range("N4").Value = 0.079926
$ws.Range("O4").Value = 0.001297498167494471
$ws.Range("P4").Value = 0.001297498167494471
$ws.Range("Q4").Value = 1.011937109311333
$ws.Range("R4").Value = 9.107433983801998
$ws.Range("S4").Value = 0.0009054599336049254
$ws.Range("T4").Value = 0.0009054599336049255
$ws.Range("D5").Value = "MuSCs"
$ws.Range("G5").Value = 37.98277566666666
$ws.Range("H5").Value = 113.948327
$ws.Range("I5").Value = 0.697850645410475
$ws.Range("J5").Value = 0.6978506454104751
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.391608
$ws.Range("N5").Value = 1.174824
$ws.Range("O5").Value = 0.01907179124600912
$ws.Range("P5").Value = 0.01907179124600912
$ws.Range("Q5").Value = 14.874358813272
$ws.Range("R5").Value = 133.869229319448
$ws.Range("S5").Value = 0.01330926183016131
$ws.Range("T5").Value = 0.01330926183016131
$ws.Range("I6").Value = 0.1779541659542351
$ws.Range("J6").Value = 0.1779541659542352
$ws.Range("M6").Value = 8.540560666666666
$ws.Range("N6").Value = 25.621682
$ws.Range("O6").Value = 0.4159358086620884
$ws.Range("P6").Value = 0.4159358086620884
$ws.Range("Q6").Value = 82.72156751280689
$ws.Range("R6").Value = 744.4941076152621
$ws.Range("S6").Value = 0.07401750992096227
$ws.Range("T6").Value = 0.07401750992096229
$ws.Range("I7").Value = 0.1779541659542351
$ws.Range("J7").Value = 0.1779541659542352
$ws.Range("O7").Value = 0.563694901924408
$ws.Range("P7").Value = 0.563694901924408
$ws.Range("S7").Value = 0.1003118561246124
$ws.Range("T7").Value = 0.1003118561246124
$ws.Range("D8").Value = "Inflammatory-Mac"
$ws.Range("I8").Value = 0.1779541659542351
$ws.Range("J8").Value = 0.1779541659542352
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.026642
$ws.Range("N8").Value = 0.079926
$ws.Range("O8").Value = 0.001297498167494471
$ws.Range("P8").Value = 0.001297498167494471
$ws.Range("Q8").Value = 0.2580472275406667
$ws.Range("R8").Value = 2.322425047866
$ws.Range("S8").Value = 0.0002308952042236271
$ws.Range("T8").Value = 0.0002308952042236271
$ws.Range("D9").Value = "MuSCs"
$ws.Range("I9").Value = 0.1779541659542351
$ws.Range("J9").Value = 0.1779541659542352
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.391608
$ws.Range("N9").Value = 1.174824
$ws.Range("O9").Value = 0.01907179124600912
$ws.Range("P9").Value = 0.01907179124600912
$ws.Range("Q9").Value = 3.793009484376
$ws.Range("R9").Value = 34.13708535938401
$ws.Range("S9").Value = 0.003393904704436835
$ws.Range("T9").Value = 0.003393904704436836
$ws.Range("G10").Value = 0.5676613333333332
$ws.Range("H10").Value = 1.702984
$ws.Range("I10").Value = 0.01042953867610283
$ws.Range("J10").Value = 0.01042953867610283
$ws.Range("M10").Value = 8.540560666666666
$ws.Range("N10").Value = 25.621682
$ws.Range("O10").Value = 0.4159358086620884
$ws.Range("P10").Value = 0.4159358086620884
$ws.Range("Q10").Value = 4.848146055454222
$ws.Range("R10").Value = 43.63331449908799
$ws.Range("S10").Value = 0.004338018603217358
$ws.Range("T10").Value = 0.004338018603217358
$ws.Range("G11").Value = 0.5676613333333332
$ws.Range("H11").Value = 1.702984
$ws.Range("I11").Value = 0.01042953867610283
$ws.Range("J11").Value = 0.01042953867610283
$ws.Range("O11").Value = 0.563694901924408
$ws.Range("P11").Value = 0.563694901924408
$ws.Range("Q11").Value = 6.570425431835554
$ws.Range("R11").Value = 59.13382888651999
$ws.Range("S11").Value = 0.005879077781142605
$ws.Range("T11").Value = 0.005879077781142605
$ws.Range("D12").Value = "Inflammatory-Mac"
$ws.Range("G12").Value = 0.5676613333333332
$ws.Range("H12").Value = 1.702984
$ws.Range("I12").Value = 0.01042953867610283
$ws.Range("J12").Value = 0.01042953867610283
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.026642
$ws.Range("N12").Value = 0.079926
$ws.Range("O12").Value = 0.001297498167494471
$ws.Range("P12").Value = 0.001297498167494471
$ws.Range("Q12").Value = 0.01512363324266666
$ws.Range("R12").Value = 0.136112699184
$ws.Range("S12").Value = 0.00001353230732005614
$ws.Range("T12").Value = 0.00001353230732005614
$ws.Range("D13").Value = "MuSCs"
$ws.Range("G13").Value = 0.5676613333333332
$ws.Range("H13").Value = 1.702984
$ws.Range("I13").Value = 0.01042953867610283
$ws.Range("J13").Value = 0.01042953867610283
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.391608
$ws.Range("N13").Value = 1.174824
$ws.Range("O13").Value = 0.01907179124600912
$ws.Range("P13").Value = 0.01907179124600912
$ws.Range("Q13").Value = 0.222300719424
$ws.Range("R13").Value = 2.000706474816
$ws.Range("S13").Value = 0.0001989099844228115
$ws.Range("T13").Value = 0.0001989099844228115
$ws.Range("G14").Value = 5.823095333333334
$ws.Range("H14").Value = 17.469286
$ws.Range("I14").Value = 0.1069866739681064
$ws.Range("J14").Value = 0.1069866739681064
$ws.Range("M14").Value = 8.540560666666666
$ws.Range("N14").Value = 25.621682
$ws.Range("O14").Value = 0.4159358086620884
$ws.Range("P14").Value = 0.4159358086620884
$ws.Range("Q14").Value = 49.73249896211689
$ws.Range("R14").Value = 447.592490659052
$ws.Range("S14").Value = 0.04449958875299155
$ws.Range("T14").Value = 0.04449958875299155
$ws.Range("G15").Value = 5.823095333333334
$ws.Range("H15").Value = 17.469286
$ws.Range("I15").Value = 0.1069866739681064
$ws.Range("J15").Value = 0.1069866739681064
$ws.Range("O15").Value = 0.563694901924408
$ws.Range("P15").Value = 0.563694901924408
$ws.Range("Q15").Value = 67.39971779559222
$ws.Range("R15").Value = 606.59746016033
$ws.Range("S15").Value = 0.06030784268967036
$ws.Range("T15").Value = 0.06030784268967036
$ws.Range("D16").Value = "Inflammatory-Mac"
$ws.Range("G16").Value = 5.823095333333334
$ws.Range("H16").Value = 17.469286
$ws.Range("I16").Value = 0.1069866739681064
$ws.Range("J16").Value = 0.1069866739681064
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.026642
$ws.Range("N16").Value = 0.079926
$ws.Range("O16").Value = 0.001297498167494471
$ws.Range("P16").Value = 0.001297498167494471
$ws.Range("Q16").Value = 0.1551389058706667
$ws.Range("R16").Value = 1.396250152836
$ws.Range("S16").Value = 0.0001388150134199466
$ws.Range("T16").Value = 0.0001388150134199466
$ws.Range("D17").Value = "MuSCs"
$ws.Range("G17").Value = 5.823095333333334
$ws.Range("H17").Value = 17.469286
$ws.Range("I17").Value = 0.1069866739681064
$ws.Range("J17").Value = 0.1069866739681064
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.391608
$ws.Range("N17").Value = 1.174824
$ws.Range("O17").Value = 0.01907179124600912
$ws.Range("P17").Value = 0.01907179124600912
$ws.Range("Q17").Value = 2.280370717296
$ws.Range("R17").Value = 20.523336455664
$ws.Range("S17").Value = 0.002040427512024564
$ws.Range("T17").Value = 0.002040427512024564
$ws.Range("G18").Value = 0.3689676666666666
$ws.Range("H18").Value = 1.106903
$ws.Range("I18").Value = 0.006778975991080511
$ws.Range("J18").Value = 0.006778975991080512
$ws.Range("M18").Value = 8.540560666666666
$ws.Range("N18").Value = 25.621682
$ws.Range("O18").Value = 0.4159358086620884
$ws.Range("P18").Value = 0.4159358086620884
$ws.Range("Q18").Value = 3.151190741205111
$ws.Range("R18").Value = 28.360716670846
$ws.Range("S18").Value = 0.002819618860750954
$ws.Range("T18").Value = 0.002819618860750955
$ws.Range("G19").Value = 0.3689676666666666
$ws.Range("H19").Value = 1.106903
$ws.Range("I19").Value = 0.006778975991080511
$ws.Range("J19").Value = 0.006778975991080512
$ws.Range("O19").Value = 0.563694901924408
$ws.Range("P19").Value = 0.563694901924408
$ws.Range("Q19").Value = 4.270635321162778
$ws.Range("R19").Value = 38.435717890465
$ws.Range("S19").Value = 0.003821274206440045
$ws.Range("T19").Value = 0.003821274206440045
$ws.Range("D20").Value = "Inflammatory-Mac"
$ws.Range("G20").Value = 0.3689676666666666
$ws.Range("H20").Value = 1.106903
$ws.Range("I20").Value = 0.006778975991080511
$ws.Range("J20").Value = 0.006778975991080512
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = 0.3333333333333333
$ws.Range("M20").Value = 0.026642
$ws.Range("N20").Value = 0.079926
$ws.Range("O20").Value = 0.001297498167494471
$ws.Range("P20").Value = 0.001297498167494471
$ws.Range("Q20").Value = 0.009830036575333332
$ws.Range("R20").Value = 0.08847032917799999
$ws.Range("S20").Value = 0.00000879570892591598
$ws.Range("T20").Value = 0.00000879570892591598
$ws.Range("D21").Value = "MuSCs"
$ws.Range("G21").Value = 0.3689676666666666
$ws.Range("H21").Value = 1.106903
$ws.Range("I21").Value = 0.006778975991080511
$ws.Range("J21").Value = 0.006778975991080512
$ws.Range("K21").Value = 3
$ws.Range("L21").Value = 1
$ws.Range("M21").Value = 0.391608
$ws.Range("N21").Value = 1.174824
$ws.Range("O21").Value = 0.01907179124600912
$ws.Range("P21").Value = 0.01907179124600912
$ws.Range("Q21").Value = 0.144490690008
$ws.Range("R21").Value = 1.300416210072
$ws.Range("S21").Value = 0.0001292872149635952
$ws.Range("T21").Value = 0.0001292872149635953

Write-Output "Applied 250 cell updates"
